# Update the date heading.
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2025-02-15 Saturday"

# Update the practice-problem table. Addressing cells directly via
# Table.Cell(row, col) avoids ambiguity from duplicate text values
# (e.g. "130×2=260" occurs twice in the original table) that a
# document-wide Find/Replace could not disambiguate.
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("388×3=1164", "881×7=6167", "826×3=2478", "506×2=1012", "341×5=1705")
    5  = @("517×6=3102", "784×5=3920", "219×6=1314", "238×2=476",  "856×8=6848")
    10 = @("988×2=1976", "772×8=6176", "646×7=4522", "486×2=972",  "454×8=3632")
    15 = @("339×2=678",  "275×3=825",  "854×8=6832", "580×7=4060", "586×8=4688")
    20 = @("124×8=992",  "879×3=2637", "536×3=1608", "128×5=640",  "433×7=3031")
}

foreach ($rowIndex in $newValues.Keys) {
    $values = $newValues[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]
    }
}
